$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (F1:H1) - copy the header style from an existing header cell (E1)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the style of the existing header row into the new header cells
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats

# Boolean data for F2:H6
$ws.Range("F2").Value = $false
$ws.Range("G2").Value = $false
$ws.Range("H2").Value = $true

$ws.Range("F3").Value = $false
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = $false

$ws.Range("F4").Value = $false
$ws.Range("G4").Value = $true
$ws.Range("H4").Value = $false

$ws.Range("F5").Value = $false
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = $false

$ws.Range("F6").Value = $false
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = $false
